# Rebuild the sensor-reading rows (2..31) of Sheet1.
# The "ax..gz" sample window was re-collected / re-sampled: the 20 existing
# sample rows shift by 3 (3 brand-new samples land on top), and 10 further
# new samples are appended at the bottom, extending the sheet from 21 to 31 rows.
# timestamp/label columns keep their original per-row values for rows 2-21,
# and rows 22-31 get new timestamps (2000..2900) with label "falling".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 30,8
$data[0,0] = 0
$data[0,1] = "falling"
$data[0,2] = -0.084752082824707
$data[0,3] = 0.6204710006713867
$data[0,4] = -1.197814345359802
$data[0,5] = -0.03934990153426188
$data[0,6] = 0.01354811775187661
$data[0,7] = -0.07257660726706182
$data[1,0] = 100
$data[1,1] = "falling"
$data[1,2] = -0.5587072372436523
$data[1,3] = 0.5920883417129517
$data[1,4] = -0.9495211839675904
$data[1,5] = -0.02474731490725574
$data[1,6] = -0.0707512873091867
$data[1,7] = 0.004276057793980506
$data[2,0] = 200
$data[2,1] = "falling"
$data[2,2] = -0.2052898406982422
$data[2,3] = 0.80674147605896
$data[2,4] = -1.046440482139587
$data[2,5] = 0.003250675749898542
$data[2,6] = 0.03200497691120401
$data[2,7] = -0.08288132186446855
$data[3,0] = 300
$data[3,1] = "falling"
$data[3,2] = -0.2705469131469726
$data[3,3] = 0.4402385950088501
$data[3,4] = -1.731513738632202
$data[3,5] = -0.05105815259351689
$data[3,6] = -0.02509637922048548
$data[3,7] = -0.0663297846913338
$data[4,0] = 400
$data[4,1] = "falling"
$data[4,2] = -0.6277971267700195
$data[4,3] = 0.0076048374176025
$data[4,4] = -1.507715225219727
$data[4,5] = -0.1197514058578581
$data[4,6] = -0.3090105539276482
$data[4,7] = 0.06041020864532091
$data[5,0] = 500
$data[5,1] = "falling"
$data[5,2] = -2.914698600769043
$data[5,3] = -1.449564576148987
$data[5,4] = -3.32840347290039
$data[5,5] = -0.4005748778581615
$data[5,6] = -0.7805985297475537
$data[5,7] = 0.1577123148100716
$data[6,0] = 600
$data[6,1] = "falling"
$data[6,2] = 1.028462886810303
$data[6,3] = -0.5502710342407227
$data[6,4] = -4.842555046081543
$data[6,5] = -0.6522004490806942
$data[6,6] = -1.329693669364566
$data[6,7] = 0.1308778794038864
$data[7,0] = 700
$data[7,1] = "falling"
$data[7,2] = -1.341280460357666
$data[7,3] = -2.225003957748413
$data[7,4] = -6.344600677490234
$data[7,5] = -0.6102398293358944
$data[7,6] = -1.220246967815219
$data[7,7] = 0.4809618578070668
$data[8,0] = 800
$data[8,1] = "falling"
$data[8,2] = 8.574896812438965
$data[8,3] = 0.6133027076721191
$data[8,4] = -6.888121604919434
$data[8,5] = -0.04692753723689685
$data[8,6] = -0.2711587122508481
$data[8,7] = 1.261436768940515
$data[9,0] = 900
$data[9,1] = "falling"
$data[9,2] = -6.096681118011475
$data[9,3] = 0.8472604751586914
$data[9,4] = 14.72706890106201
$data[9,5] = 0.2368920927955998
$data[9,6] = 2.034886604263666
$data[9,7] = 0.6732607796078658
$data[10,0] = 1000
$data[10,1] = "falling"
$data[10,2] = 4.274323463439941
$data[10,3] = -4.468049049377441
$data[10,4] = -6.856836795806885
$data[10,5] = -0.7203119397163371
$data[10,6] = 3.954537868499759
$data[10,7] = -2.220546166102086
$data[11,0] = 1100
$data[11,1] = "falling"
$data[11,2] = -4.518700122833252
$data[11,3] = -1.648021101951599
$data[11,4] = -0.9248533248901368
$data[11,5] = -0.9733701603753284
$data[11,6] = 2.440581185477122
$data[11,7] = -3.891320841653013
$data[12,0] = 1200
$data[12,1] = "falling"
$data[12,2] = 9.755411148071287
$data[12,3] = 3.367114305496216
$data[12,4] = 2.822277307510376
$data[12,5] = 0.4988514525549754
$data[12,6] = 0.8867653551555867
$data[12,7] = -1.71078631139937
$data[13,0] = 1300
$data[13,1] = "falling"
$data[13,2] = 1.561064720153809
$data[13,3] = 0.1129603385925293
$data[13,4] = -0.9029455184936525
$data[13,5] = 0.1264272814705268
$data[13,6] = -2.447191684019
$data[13,7] = -0.4622068021978656
$data[14,0] = 1400
$data[14,1] = "falling"
$data[14,2] = 5.92741584777832
$data[14,3] = -0.8555939197540283
$data[14,4] = 4.797466278076172
$data[14,5] = -0.7295694393771054
$data[14,6] = -0.004014266388761123
$data[14,7] = 0.09374600010258789
$data[15,0] = 1500
$data[15,1] = "falling"
$data[15,2] = 1.122594833374023
$data[15,3] = 1.295500755310059
$data[15,4] = -1.442571401596069
$data[15,5] = -0.1179624412740978
$data[15,6] = 1.259269575277969
$data[15,7] = 0.1631955632141657
$data[16,0] = 1600
$data[16,1] = "falling"
$data[16,2] = 0.5986118316650391
$data[16,3] = 0.4096674025058746
$data[16,4] = -0.6679027080535889
$data[16,5] = -0.1617193005624273
$data[16,6] = -0.4078179995218952
$data[16,7] = 0.2281727109636559
$data[17,0] = 1700
$data[17,1] = "falling"
$data[17,2] = 0.0388402938842773
$data[17,3] = 0.3524296283721924
$data[17,4] = -1.101761341094971
$data[17,5] = -0.124921940267086
$data[17,6] = 1.249372124671936
$data[17,7] = 1.016479730606079
$data[18,0] = 1800
$data[18,1] = "falling"
$data[18,2] = -0.1728830337524414
$data[18,3] = 0.6193998456001282
$data[18,4] = -0.6873818635940552
$data[18,5] = -0.4928955077415405
$data[18,6] = -0.4476696934018742
$data[18,7] = -0.9833766732896989
$data[19,0] = 1900
$data[19,1] = "falling"
$data[19,2] = 0.4876585006713867
$data[19,3] = 0.6636635065078735
$data[19,4] = -0.9166454076766968
$data[19,5] = -0.2805471434479678
$data[19,6] = 0.4230750912711692
$data[19,7] = -0.2188279224293611
$data[20,0] = 2000
$data[20,1] = "falling"
$data[20,2] = -0.1092472076416015
$data[20,3] = 0.732629120349884
$data[20,4] = -1.016466021537781
$data[20,5] = -0.2866266923291341
$data[20,6] = 0.4182899764605935
$data[20,7] = -0.003992439912898826
$data[21,0] = 2100
$data[21,1] = "falling"
$data[21,2] = 0.4153709411621094
$data[21,3] = 0.5096800327301025
$data[21,4] = -0.7671611309051514
$data[21,5] = -0.03713915026968551
$data[21,6] = 0.07855436143775912
$data[21,7] = 0.07685266648020034
$data[22,0] = 2200
$data[22,1] = "falling"
$data[22,2] = 0.17730712890625
$data[22,3] = 0.6253083348274231
$data[22,4] = -0.8837988376617432
$data[22,5] = 0.03008511281084441
$data[22,6] = 0.0510363349070151
$data[22,7] = -0.07897615255344478
$data[23,0] = 2300
$data[23,1] = "falling"
$data[23,2] = 0.0039987564086914
$data[23,3] = 0.546174168586731
$data[23,4] = -0.7374091148376465
$data[23,5] = 0.0188495556690863
$data[23,6] = -0.09346238630158486
$data[23,7] = -0.02356194624943445
$data[24,0] = 2400
$data[24,1] = "falling"
$data[24,2] = 0.0831842422485351
$data[24,3] = 0.5668889284133911
$data[24,4] = -0.8130950927734375
$data[24,5] = 0.1422297873844697
$data[24,6] = -0.01903863499562003
$data[24,7] = -0.03227404815455254
$data[25,0] = 2500
$data[25,1] = "falling"
$data[25,2] = -0.0142126083374023
$data[25,3] = 0.51572585105896
$data[25,4] = -0.7260744571685791
$data[25,5] = 0.08552113210871153
$data[25,6] = 0.01936588267562919
$data[25,7] = -0.03527019580914864
$data[26,0] = 2600
$data[26,1] = "falling"
$data[26,2] = 0.1546173095703125
$data[26,3] = 0.5381616353988647
$data[26,4] = -0.7814648151397705
$data[26,5] = -0.0006981316421711559
$data[26,6] = 0.01939497157166284
$data[26,7] = -0.01994038639324046
$data[27,0] = 2700
$data[27,1] = "falling"
$data[27,2] = 0.2052326202392578
$data[27,3] = 0.5754936933517456
$data[27,4] = -0.8383152484893799
$data[27,5] = 0.006530440013323374
$data[27,6] = 0.04903648190555104
$data[27,7] = -0.02060942954960315
$data[28,0] = 2800
$data[28,1] = "falling"
$data[28,2] = -0.07891082763671869
$data[28,3] = 0.5558477640151978
$data[28,4] = -0.7180624008178711
$data[28,5] = 0.00994837645529994
$data[28,6] = 0.04445499217226397
$data[28,7] = -0.01362811268440311
$data[29,0] = 2900
$data[29,1] = "falling"
$data[29,2] = 0.1948976516723632
$data[29,3] = 0.6977589726448059
$data[29,4] = -0.9572491645812988
$data[29,5] = 0.01760600972920647
$data[29,6] = -0.0257654253925594
$data[29,7] = -0.005214171284543525

$ws.Range("A2:H31").Value = $data
